$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at row 640 (shifts existing rows 640+ down by 9)
$ws.Range("A640:I648").EntireRow.Insert()

# Columns B (date) and C (id) must stay text, not auto-converted to date/number
$ws.Range("B640:C648").NumberFormat = "@"

# Populate the 9 new rows with historical data (2019-11-18 .. 2019-11-28)
$ws.Range("A640").Value2 = 1574035200
$ws.Range("B640").Value = "2019-11-18"
$ws.Range("C640").Value = "5283"
$ws.Range("D640").Value = "EWINT"
$ws.Range("E640").Value2 = 0.78
$ws.Range("F640").Value2 = 0.82
$ws.Range("G640").Value2 = 0.78
$ws.Range("H640").Value2 = 0.8149999999999999
$ws.Range("I640").Value2 = 1901400

$ws.Range("A641").Value2 = 1574121600
$ws.Range("B641").Value = "2019-11-19"
$ws.Range("C641").Value = "5283"
$ws.Range("D641").Value = "EWINT"
$ws.Range("E641").Value2 = 0.8149999999999999
$ws.Range("F641").Value2 = 0.83
$ws.Range("G641").Value2 = 0.8149999999999999
$ws.Range("H641").Value2 = 0.825
$ws.Range("I641").Value2 = 2214600

$ws.Range("A642").Value2 = 1574208000
$ws.Range("B642").Value = "2019-11-20"
$ws.Range("C642").Value = "5283"
$ws.Range("D642").Value = "EWINT"
$ws.Range("E642").Value2 = 0.83
$ws.Range("F642").Value2 = 0.83
$ws.Range("G642").Value2 = 0.8149999999999999
$ws.Range("H642").Value2 = 0.825
$ws.Range("I642").Value2 = 1948900

$ws.Range("A643").Value2 = 1574294400
$ws.Range("B643").Value = "2019-11-21"
$ws.Range("C643").Value = "5283"
$ws.Range("D643").Value = "EWINT"
$ws.Range("E643").Value2 = 0.825
$ws.Range("F643").Value2 = 0.825
$ws.Range("G643").Value2 = 0.8149999999999999
$ws.Range("H643").Value2 = 0.82
$ws.Range("I643").Value2 = 766900

$ws.Range("A644").Value2 = 1574380800
$ws.Range("B644").Value = "2019-11-22"
$ws.Range("C644").Value = "5283"
$ws.Range("D644").Value = "EWINT"
$ws.Range("E644").Value2 = 0.82
$ws.Range("F644").Value2 = 0.825
$ws.Range("G644").Value2 = 0.8100000000000001
$ws.Range("H644").Value2 = 0.825
$ws.Range("I644").Value2 = 796900

$ws.Range("A645").Value2 = 1574640000
$ws.Range("B645").Value = "2019-11-25"
$ws.Range("C645").Value = "5283"
$ws.Range("D645").Value = "EWINT"
$ws.Range("E645").Value2 = 0.825
$ws.Range("F645").Value2 = 0.945
$ws.Range("G645").Value2 = 0.825
$ws.Range("H645").Value2 = 0.91
$ws.Range("I645").Value2 = 8688500

$ws.Range("A646").Value2 = 1574726400
$ws.Range("B646").Value = "2019-11-26"
$ws.Range("C646").Value = "5283"
$ws.Range("D646").Value = "EWINT"
$ws.Range("E646").Value2 = 0.91
$ws.Range("F646").Value2 = 0.93
$ws.Range("G646").Value2 = 0.9
$ws.Range("H646").Value2 = 0.9
$ws.Range("I646").Value2 = 4694200

$ws.Range("A647").Value2 = 1574812800
$ws.Range("B647").Value = "2019-11-27"
$ws.Range("C647").Value = "5283"
$ws.Range("D647").Value = "EWINT"
$ws.Range("E647").Value2 = 0.905
$ws.Range("F647").Value2 = 0.915
$ws.Range("G647").Value2 = 0.9
$ws.Range("H647").Value2 = 0.905
$ws.Range("I647").Value2 = 963500

$ws.Range("A648").Value2 = 1574899200
$ws.Range("B648").Value = "2019-11-28"
$ws.Range("C648").Value = "5283"
$ws.Range("D648").Value = "EWINT"
$ws.Range("E648").Value2 = 0.905
$ws.Range("F648").Value2 = 0.93
$ws.Range("G648").Value2 = 0.905
$ws.Range("H648").Value2 = 0.93
$ws.Range("I648").Value2 = 3301800
